# Updated password service tests
#
# Mark the four PasswordServiceTest.java rows (67-70) and the four
# UserRepositoryTest.java rows (79-82) as "Y" in the "JS Implemented"
# column (D). D107/E107 hold COUNTIF/ROWS formulas over this column and
# will pick up the new values automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 67,68,69,70,79,80,81,82
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "Y"
}

# Restore the on-screen selection/scroll position to where the author
# left it: viewport scrolled to row 79, with D71 the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 79
$win.ScrollColumn = 1
[void]$ws.Range("D71").Select()
